$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear all cell contents (keeps formatting) so the shared-string table is
# rebuilt fresh, following first-use order as we rewrite the data below.
$ws.Cells.ClearContents()

# Header row
$ws.Range("A1").Value = "Cluster Name"
$ws.Range("B1").Value = "Active cases"

# Full data set (Cluster Name, Active cases), sorted by Cluster Name
$data = @(
    @("3433 HammondCare Caulfield Village Aged Care", 11),
    @("3642 Fronditha Care Aged Care Clayton South", 41),
    @("4314 Estia Health Altona Meadows", 13),
    @("Beyond the City New Years Festival 30 to 31 Jan Melbourne", 43),
    @("Confirmed Omicron Sircuit Bar Fitzroy", 34),
    @("Confirmed Omicron Variant The Peel Hotel Collingwood", 25),
    @("Diamond Valley Pork and Baxters Pork Laverton North", 34),
    @("Mercure Welcome Melbourne", 11),
    @("Novotel ibis Melbourne Central Melbourne", 18),
    @("Pullman Melbourne on Swanston Melbourne", 13),
    @("Werribee Mercy Hospital Emergency Department", 16),
    @("Western Health Sunshine Hospital EmergencyDepartment St Albans", 10)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row = $row + 1
}
